$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J1").Value = "Q8"
$ws.Range("J1").Font.Bold = $true
$ws.Range("J1").HorizontalAlignment = -4108
$ws.Range("J1").VerticalAlignment = -4160
$ws.Range("J1").Borders.LineStyle = 1

$ws.Range("A2").Value = "2022-01-01 00:00:00_diff"
$ws.Range("B2").Value = -0.740857461610962
$ws.Range("C2").Value = 0.2483496536967165
$ws.Range("D2").Value = 0.8032852183307098
$ws.Range("E2").Value = 1.015296315185831
$ws.Range("F2").Value = 0.4272448182495295
$ws.Range("G2").Value = 0.01108471313272752
$ws.Range("H2").Value = 1.401227768176947

$ws.Range("A3").Value = "2022-04-01 00:00:00_diff"
$ws.Range("B3").Value = 0.2917404373296907
$ws.Range("C3").Value = 0.8466760019636841
$ws.Range("D3").Value = 1.058687098818806
$ws.Range("E3").Value = 0.4706356018825037
$ws.Range("F3").Value = 0.05447549676570174
$ws.Range("G3").Value = 1.444618551809921

$ws.Range("A4").Value = "2022-07-01 00:00:00_diff"
$ws.Range("B4").Value = 1.247683574918372
$ws.Range("C4").Value = 1.459694671773494
$ws.Range("D4").Value = 0.8716431748371921
$ws.Range("E4").Value = 0.4554830697203902
$ws.Range("F4").Value = 1.84562612476461
$ws.Range("G4").Value = 1.076491995083831
$ws.Range("H4").Value = -0.7442442573846902
$ws.Range("I4").Value = 1.069262038377534
$ws.Range("J4").Value = 0.4350014876132097

$ws.Range("A5").Value = "2022-10-01 00:00:00_diff"
$ws.Range("B5").Value = 3.694617372890321
$ws.Range("C5").Value = 3.106565875954019
$ws.Range("D5").Value = 2.690405770837217
$ws.Range("E5").Value = 4.080548825881436
$ws.Range("F5").Value = 3.311414696200658
$ws.Range("G5").Value = 1.490678443732137
$ws.Range("H5").Value = 3.304184739494361
$ws.Range("I5").Value = 2.669924188730037

$ws.Range("A6").Value = "2023-01-01 00:00:00_diff"
$ws.Range("B6").Value = 1.216530487278416
$ws.Range("C6").Value = 0.8003703821616144
$ws.Range("D6").Value = 2.190513437205834
$ws.Range("E6").Value = 1.421379307525055
$ws.Range("F6").Value = -0.399356944943466
$ws.Range("G6").Value = 1.414149350818758
$ws.Range("H6").Value = 0.779888800054434

$ws.Range("A7").Value = "2023-04-01 00:00:00_diff"
$ws.Range("B7").Value = 0.4476790584865185
$ws.Range("C7").Value = 1.837822113530738
$ws.Range("D7").Value = 1.068687983849959
$ws.Range("E7").Value = -0.7520482686185619
$ws.Range("F7").Value = 1.061458027143662
$ws.Range("G7").Value = 0.427197476379338

$ws.Range("A8").Value = "2023-07-01 00:00:00_diff"
$ws.Range("B8").Value = 1.594404170131267
$ws.Range("C8").Value = 0.8252700404504878
$ws.Range("D8").Value = -0.9954662120180333
$ws.Range("E8").Value = 0.8180400837441908
$ws.Range("F8").Value = 0.1837795329798666
$ws.Range("G8").Value = 1.265495818666463
$ws.Range("H8").Value = 0.8840541853673727
$ws.Range("I8").Value = 0.5915919440004813

$ws.Range("A9").Value = "2023-10-01 00:00:00_diff"
$ws.Range("B9").Value = 0.4048306212132332
$ws.Range("C9").Value = -1.415905631255288
$ws.Range("D9").Value = 0.3976006645069362
$ws.Range("E9").Value = -0.236659886257388
$ws.Range("F9").Value = 0.8450563994292083
$ws.Range("G9").Value = 0.4636147661301181
$ws.Range("H9").Value = 0.1711525247632267

$ws.Range("A10").Value = "2024-01-01 00:00:00_diff"
$ws.Range("B10").Value = -1.433992460878194
$ws.Range("C10").Value = 0.3795138348840296
$ws.Range("D10").Value = -0.2547467158802946
$ws.Range("E10").Value = 0.8269695698063018
$ws.Range("F10").Value = 0.4455279365072115
$ws.Range("G10").Value = 0.1530656951403201

$ws.Range("A11").Value = "2024-04-01 00:00:00_diff"
$ws.Range("B11").Value = 0.4299722955860048
$ws.Range("C11").Value = -0.2042882551783194
$ws.Range("D11").Value = 0.8774280305082769
$ws.Range("E11").Value = 0.4959863972091867
$ws.Range("F11").Value = 0.2035241558422953

$ws.Range("A12").Value = "2024-07-01 00:00:00_diff"
$ws.Range("B12").Value = -0.4160968922281114
$ws.Range("C12").Value = 0.6656193934584849
$ws.Range("D12").Value = 0.2841777601593947
$ws.Range("E12").Value = -0.008284481207496679

$ws.Range("A13").Value = "2024-10-01 00:00:00_diff"
$ws.Range("B13").Value = 0.5354267536149976
$ws.Range("C13").Value = 0.1539851203159074
$ws.Range("D13").Value = -0.1384771210509839

$ws.Range("A14").Value = "2025-01-01 00:00:00_diff"
$ws.Range("B14").Value = -0.322788625881465
$ws.Range("C14").Value = -0.6152508672483563

$ws.Range("A15").Value = "2025-04-01 00:00:00_diff"
$ws.Range("B15").Value = -0.2720993704486361

$ws.Range("A16").Value = "2025-07-01 00:00:00_diff"
